$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 13 (shifts old rows 13-23 down to 15-25)
$ws.Rows.Item(13).Insert()
$ws.Rows.Item(13).Insert()

# Update cell contents (B and C columns mirror each other throughout this sheet)
$ws.Range("B10:C10").Value = 'Apresentar uma visão técnico-científica e mercadológica da indústria de polímeros termoplásticos, popularmente conhecidos como plásticos e dos elastômeros.'
$ws.Range("B13:C13").Value = '5840897 - Clodoaldo Saron'
$ws.Range("B14:C14").Value = '1033242 - Fábio Herbst Florenzano'
$ws.Range("B15:C15").Value = 'Principais tipos de polímeros termoplásticos ("commodities" e de engenharia) e elastômeros. Mercados, tecnologias de trasnformação e reciclagem.'
$ws.Range("B17:C17").Value = '1. Definição e classificação de polímeros termoplásticos, elastômeros e fibras.2. Identificação de plásticos, borrachas e fibras.3. Polímeros de adição olefínicos: polietileno, polipropileno e seus copolímeros.4. Polímeros de condensação: poli(tereftalato de etileno), poli(tereftalato de butileno) e poliamidas.5. Polímeros halogenados: poli(cloreto de vinila), poli(tetrafluor etileno) e poli(fluoreto de vinilideno).6. Termoplásticos acrílicos e oximetilênicos: PMMA, POM e poliacetais.7. Termoplásticos nitrogenados: poliacrilonitrila, poliuretano, ABS e SAN.8. Termoplásticos estirênicos e fenólicos: poliestireno, HIPS, SBR e policarbonato.9. Polímeros hidrolisáveis: EVA, PVAc e PEO.10. Termoplásticos avançados: PPO, PPS e PEEK.11. Elastômeros: borracha natural, polibutadieno, borrachas nitrílicas e fluoradas, EPDM e polisiloxanos.12. Aditivos e compostos.13. Tecnologias de transformação apropriadas a cada tipo de plástico: extrusão, injeção, laminação, calandragem, termoformação e moldagem por sopro.14. Testes e ensaios de polímeros termoplásticos e elastômeros.15. Reciclagem.'
$ws.Range("B20:C20").Value = 'A avaliação será feita por meio de Provas Escritas, Estudos de Casos e Desenvolvimento de Projetos, sendo necessário utilizar pelo menos dois critérios de avaliação diferentes.'
$ws.Range("B21:C21").Value = 'A Nota final (NF) será calculada da seguinte maneira: NF = (P+EC+Projetos)/3'
$ws.Range("B22:C22").Value = 'Não consta recuperação'
$ws.Range("B23:C23").Value = '1. J. Margolis. Engineering Plastics Handbook. McGraw-Hill Professional, 2005. 2. Nigel Mills. Plastics - Microstructure and Engineering Applications. Butterworth-Heineman, 2005. 3. Walter Michaeli, TEcnologia dos Plasticos. Ed. Blucher 4. Hélio Wiebeck, Júlio Harada. Plásticos de Engenharia - Tecnologia e Aplicações. São Paulo: Editora Artliber, 2005. 5. E. B. Mano, L. C. Mendes. Identificação de Plásticos, Borrachas e Fibras. São Paulo: Editora Edgard Blucher, 2000. 6. Marcelo Rabello. Aditivação de Polímeros. São Paulo: Editora Artliber, 2004. 7. Jan C.J. Bart. Additives in Polymers. New York: John Wiley & Sons, 2005. 8. Marino Xanthos. Functional Fillers for Plastics. Wiley-VCH Verlag GmbH, 2005. 9. Silvio Manrich. Processamento de Termoplásticos. Editora Artliber, 2005. 10. G.H. Michler, F.J. Baltá-Calleja. Mechanical Properties of Polymers Based on Nanostructure and Morphology. Boca Raton: CRC Press, 2005. 11. A. M. Piva, H. Wiebeck. Reciclagem do P. São Paulo: Editora Artliber". Manas Chanda, ,Salil K. Roy  Plastics Fabrication and Recycling'

Write-Output "done"